$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.502.91"
$ws.Range("E2").Value = "  -1.67%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.670.75"
$ws.Range("E3").Value = "  -1.93%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.49"
$ws.Range("E5").Value = "  -1.01%  "

# Row 6
$ws.Range("E6").Value = "  +0.00%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3904"
$ws.Range("E7").Value = "  -3.88%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3930"
$ws.Range("E8").Value = "  -3.35%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.002"
$ws.Range("E9").Value = "  -0.01%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "51.62"
$ws.Range("E10").Value = "  -3.98%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.396"
$ws.Range("E11").Value = "  -4.92%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08626"
$ws.Range("E12").Value = "  -2.22%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.33"
$ws.Range("E13").Value = "  -1.95%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.288"
$ws.Range("E14").Value = "  -3.18%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001318"
$ws.Range("E15").Value = "  -2.51%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.699"
$ws.Range("E16").Value = "  -4.49%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.672.25"
$ws.Range("E17").Value = "  -1.72%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.24"
$ws.Range("E18").Value = "  -3.48%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07035"
$ws.Range("E19").Value = "  -2.04%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.86"
$ws.Range("E20").Value = "  -0.89%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.052"
$ws.Range("E21").Value = "  -2.67%  "

# Row 22
$ws.Range("E22").Value = "  +0.08%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.93"
$ws.Range("E23").Value = "  -4.73%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.507.87"
$ws.Range("E24").Value = "  -1.62%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.368"
$ws.Range("E25").Value = "  +1.76%  "

# Row 26
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.739"
$ws.Range("E26").Value = "  -5.27%  "

# Row 27
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.13"
$ws.Range("E27").Value = "  +0.29%  "

# Row 28
$ws.Range("B28").Value = "HuobiToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.885"
$ws.Range("E28").Value = "  -13.35%  "

# Row 29
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "160.72"
$ws.Range("E29").Value = "  -2.62%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "147.28"
$ws.Range("E30").Value = "  +1.39%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.254"
$ws.Range("E31").Value = "  -0.32%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.508"
$ws.Range("E32").Value = "  +10.43%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.856.97"
$ws.Range("E33").Value = "  -2.57%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08347"
$ws.Range("E34").Value = "  -5.00%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.983"
$ws.Range("E35").Value = "  -4.70%  "

# Row 36
$ws.Range("E36").Value = "  -5.86%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2805"
$ws.Range("E37").Value = "  -1.58%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9843"
$ws.Range("E38").Value = "  -3.21%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09451"
$ws.Range("E39").Value = "  +0.44%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.521"
$ws.Range("E40").Value = "  +3.59%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.31"
$ws.Range("E41").Value = "  -4.59%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7888"
$ws.Range("E42").Value = "  -7.15%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.58"
$ws.Range("E43").Value = "  -3.42%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.40"
$ws.Range("E44").Value = "  -8.38%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7111"
$ws.Range("E45").Value = "  -4.55%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.548"
$ws.Range("E46").Value = "  -6.26%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.173"
$ws.Range("E47").Value = "  -1.61%  "

# Row 48
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.08607"
$ws.Range("E48").Value = "  +2.95%  "

# Row 49
$ws.Range("B49").Value = "Frax"
$ws.Range("C49").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.002"
$ws.Range("E49").Value = "  -0.05%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.318"
$ws.Range("E50").Value = "  -5.30%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "137.47"
$ws.Range("E51").Value = "  -3.22%  "
